$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the season record columns
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Copy style from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in season record values (Wins, Losses, Ties) for each data row
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 90
    $ws.Cells.Item($r, 31).Value2 = 73
    $ws.Cells.Item($r, 32).Value2 = 0
}
